# Generate Report for Handback
#
# For the "40fa5254-b5ae-4d84-9c75-7c8701305610" row (row 6) on both the
# zh-cn and de-de sheets, the handback pipeline discovered that the file
# handed back was not built from the latest source: it now fills in the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# / "Error Detail" columns (I, J, K, P) that were previously blank, and
# widens column P (Error Detail) so the message is readable.

$wb = $excel.ActiveWorkbook

$targetFileName = "40fa5254-b5ae-4d84-9c75-7c8701305610.md"
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5ca49cd1e328e588fdca2f351bb39676ae8a1915/e2e/40fa5254-b5ae-4d84-9c75-7c8701305610.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d6e644b1ecc4345d30e216f7da2ca7a675fa771e/e2e/40fa5254-b5ae-4d84-9c75-7c8701305610.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5ca49cd1e328e588fdca2f351bb39676ae8a1915/e2e/40fa5254-b5ae-4d84-9c75-7c8701305610.md."

function Update-HandbackSheet($sheetName, $xlfFileName, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Error Detail column (P) is too narrow for the new message - widen it.
    $ws.Columns.Item(16).ColumnWidth = 39.15

    # Latest Target File
    $ws.Range("I6").Value = $targetFileName
    $ws.Hyperlinks.Add($ws.Range("I6"), $targetUrl, [Type]::Missing, [Type]::Missing, $targetFileName) | Out-Null
    $ws.Range("I6").Font.Underline = 2
    $ws.Range("I6").Font.Color = 15570276

    # Latest Handback File
    $ws.Range("J6").Value = $xlfFileName

    # Latest Handback DateTime
    $ws.Range("K6").Value = $handbackDateTime

    # Error Detail
    $ws.Range("P6").Value = $errorDetail
}

Update-HandbackSheet "zh-cn" "40fa5254-b5ae-4d84-9c75-7c8701305610.adee70e4b479b933c4feed746aaff4a94a835711.zh-cn.xlf" "2016-09-06 15:29:30"
Update-HandbackSheet "de-de" "40fa5254-b5ae-4d84-9c75-7c8701305610.adee70e4b479b933c4feed746aaff4a94a835711.de-de.xlf" "2016-09-06 15:29:39"
